$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the formatting
# (bold font, border, centered alignment) already used by the other
# header cells (e.g. H1) by copying H1's formatting onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I / J data columns for rows 2-13.
$values = @(
  @(6, 8),
  @(10, 10),
  @(4, 5),
  @(9, 9),
  @(10, 10),
  @(7, 8),
  @(6, 8),
  @(6, 7),
  @(5, 7),
  @(7, 7),
  @(4, 5),
  @(3, 4)
)

$rowCount = $values.Count
$data = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
  $data[$i, 0] = $values[$i][0]
  $data[$i, 1] = $values[$i][1]
}

$ws.Range("I2:J13").Value = $data
